$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("100_1")

# Fix typo'd text values that should be numeric entries.
# B14 held the text "ll0" (letters L, L, zero) -> correct numeric value 110
$ws.Range("B14").Value = 110
# B40 held the text "l" (letter L) -> correct numeric value 1
$ws.Range("B40").Value = 1
